{"js": "// The diary has two \"21/04/2015\" entry tables. In the first one the date is\n// split across two runs (\"21\" + \"/04/2015\"); `search` matches across the run\n// boundary, and re-inserting the same text via `insertText(..., \"Replace\")`\n// collapses the match back down into a single run. The second table's date\n// cell is already a single run with the same text, so re-writing it is a\n// harmless no-op.\nconst dateHits = context.document.body.search(\"21/04/2015\", { matchCase: true });\nawait context.sync();\ndateHits.items.forEach((r) => r.insertText(\"21/04/2015\", \"Replace\"));\nawait context.sync();\n\n// \"Time Taken\" for the final (Report) entry: \"3 hours 50 minutes.\" -> \"5 hours.\"\nconst timeTaken = context.document.body.search(\"3 hours 50 minutes.\", { matchCase: true });\nawait context.sync();\ntimeTaken.items.forEach((r) => r.insertText(\"5 hours.\", \"Replace\"));\nawait context.sync();\n\n// \"Problems\" cell: append a sentence about reformatting / contents page.\nconst problems = context.document.body.search(\n  \"Everyone\\u2019s document had a completely different format.\",\n  { matchCase: true }\n);\nawait context.sync();\nproblems.items.forEach((r) =>\n  r.insertText(\n    \"Everyone\\u2019s document had a completely different format. Reformatting and making a contents page.\",\n    \"Replace\"\n  )\n);\nawait context.sync();\n\n// \"Solutions\" cell: append a sentence about the process being lengthy/delicate.\nconst solutions = context.document.body.search(\n  \"I had to go through each document and make it as consistent as possible.\",\n  { matchCase: true }\n);\nawait context.sync();\nsolutions.items.forEach((r) =>\n  r.insertText(\n    \"I had to go through each document and make it as consistent as possible. This was an extremely lengthy and delicate process.\",\n    \"Replace\"\n  )\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Constants (kept literal since the COM enum names aren't bound in this host):\n#   wdFindContinue = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll) | Out-Null\n}\n\n# 1. The date \"21/04/2015\" in the first diary entry is split across two runs\n#    (\"21\" + \"/04/2015\"). Re-running Find/Replace over the same visible text\n#    collapses the match back into a single run (the second, already-merged\n#    \"21/04/2015\" table cell is left as a harmless no-op).\nReplace-Text \"21/04/2015\" \"21/04/2015\"\n\n# 2. Final entry's \"Time Taken\": \"3 hours 50 minutes.\" -> \"5 hours.\"\nReplace-Text \"3 hours 50 minutes.\" \"5 hours.\"\n\n# 3. Final entry's \"Problems\": append a sentence about reformatting.\n$rsquo = [char]0x2019\n$problemsOld = \"Everyone${rsquo}s document had a completely different format.\"\n$problemsNew = \"$problemsOld Reformatting and making a contents page.\"\nReplace-Text $problemsOld $problemsNew\n\n# 4. Final entry's \"Solutions\": append a sentence about the lengthy process.\n$solutionsOld = \"I had to go through each document and make it as consistent as possible.\"\n$solutionsNew = \"$solutionsOld This was an extremely lengthy and delicate process.\"\nReplace-Text $solutionsOld $solutionsNew\n"}
